$wb = $excel.ActiveWorkbook

# --- optimized_production_rates ---
$ws = $wb.Worksheets.Item("optimized_production_rates")
$ws.Range("B2").Value = 0.31340589843110728
$ws.Range("B3").Value = 0.81784399356117932
$ws.Range("B4").Value = 0.31897356641545183

# --- optimized_threshold_b ---
$ws = $wb.Worksheets.Item("optimized_threshold_b")
$ws.Range("B2").Value = 0.040704848277978226
$ws.Range("B3").Value = 1.2763252158261853
$ws.Range("B4").Value = -0.068341714754867278

# --- network_optimized_weights ---
$ws = $wb.Worksheets.Item("network_optimized_weights")
$ws.Range("C2").Value = -1.9794090510770921
$ws.Range("D3").Value = -0.89020246593128949
$ws.Range("B4").Value = 1.6517034868941354
$ws.Range("C4").Value = 0.80355636273399977

# --- optimization_diagnostics ---
$ws = $wb.Worksheets.Item("optimization_diagnostics")
$ws.Range("B2").Value = 0.0000000011307916401083287
$ws.Range("B3").Value = 1.0588517580794632
$ws.Range("B5").Value = 10103
$ws.Range("B8").Value = 0.00000000016896798687566769
$ws.Range("B9").Value = 0.0000000011175099036090508
$ws.Range("B10").Value = 0.0000000021058970298402674

# --- network_weights ---
$ws = $wb.Worksheets.Item("network_weights")
$ws.Range("C2").Value = -2
$ws.Range("D3").Value = -1
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1

# --- wt_log2_optimized_expression ---
$ws = $wb.Worksheets.Item("wt_log2_optimized_expression")
$ws.Range("C2").Value = -0.6435055754054021
$ws.Range("D2").Value = -0.85894189805469767
$ws.Range("E2").Value = -0.82719445806623282
$ws.Range("F2").Value = -0.73877142795677853
$ws.Range("G2").Value = -0.66313735119876516
$ws.Range("H2").Value = -0.61077673205877914
$ws.Range("I2").Value = -0.57695688340669071
$ws.Range("J2").Value = -0.55550810136892814
$ws.Range("K2").Value = -0.5418847689720081
$ws.Range("L2").Value = -0.53316051868150793
$ws.Range("M2").Value = -0.52752434354349509
$ws.Range("N2").Value = -0.52385631081463258
$ws.Range("C3").Value = -0.51750644128539913
$ws.Range("D3").Value = -0.95863304853420639
$ws.Range("E3").Value = -1.2580780694677571
$ws.Range("F3").Value = -1.4353727778564955
$ws.Range("G3").Value = -1.5346440804020076
$ws.Range("H3").Value = -1.590988507112185
$ws.Range("I3").Value = -1.6246293397187541
$ws.Range("J3").Value = -1.6457658830162578
$ws.Range("K3").Value = -1.6594983437202535
$ws.Range("L3").Value = -1.6685663275457168
$ws.Range("M3").Value = -1.674588963940792
$ws.Range("N3").Value = -1.6785896792599782
$ws.Range("C4").Value = 0.52916817760706936
$ws.Range("D4").Value = 0.65985380951791484
$ws.Range("E4").Value = 0.69514273201567112
$ws.Range("F4").Value = 0.71032247211252897
$ws.Range("G4").Value = 0.72210675934841562
$ws.Range("H4").Value = 0.73234367857017202
$ws.Range("I4").Value = 0.74059287694435383
$ws.Range("J4").Value = 0.74672943224465416
$ws.Range("K4").Value = 0.75105527370394443
$ws.Range("L4").Value = 0.75400610900025111
$ws.Range("M4").Value = 0.75598100451980388
$ws.Range("N4").Value = 0.75728903482136767

# --- network_weights: update selected cell to match diff (C4) ---
$wsNW = $wb.Worksheets.Item("network_weights")
$wsNW.Range("C4").Select()

# --- restore original active sheet (activeTab stays unchanged in diff) ---
$wsFinal = $wb.Worksheets.Item("optimization_diagnostics")
$wsFinal.Activate()
